$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()
$ws.Select()

# Remember column B's width so the freshly inserted column C can match it.
$bWidth = $ws.Range("B1").ColumnWidth

# Insert a new column before the current column C (email). Existing C/D
# (email/password) shift right to D/E.
$ws.Columns("C").Insert()

# Give the new column the same width as column B.
$ws.Range("C1").ColumnWidth = $bWidth

# ---- Re-point the existing hyperlinks (currently still anchored to the old
#      column C) onto the new column D, preserving order/targets/rIds. ----
$links = @()
foreach ($h in $ws.Hyperlinks) {
    $links += , @($h.Range.Address(), $h.Address, $h.TextToDisplay)
}
$ws.Hyperlinks.Delete()
foreach ($l in $links) {
    $oldAddr = $l[0]
    $newAddr = $oldAddr.Replace("C", "D")
    $target = $l[1]
    $range = $ws.Range($newAddr)
    if ($range.Cells.Count -gt 1) {
        $display = $l[2]
        $ws.Hyperlinks.Add($range, $target, [Type]::Missing, [Type]::Missing, $display) | Out-Null
    } else {
        $ws.Hyperlinks.Add($range, $target) | Out-Null
    }
}

# ---- Header for the new "dateOfBirth" column ----
$ws.Range("C1").Value = "dateOfBirth"

# ---- Fill the new column with date-of-birth text values ----
# Set the number format to Text ("@") first so Excel keeps the values as
# plain strings instead of converting them to date serial numbers.
$ws.Range("C2:C22").NumberFormat = "@"
for ($r = 2; $r -le 22; $r++) {
    $dayIndex = (($r - 2) % 10) + 1
    $dateText = "2000-01-{0:D2}" -f $dayIndex
    $ws.Cells.Item($r, 3).Value = $dateText
}

# Match the page orientation recorded in the saved file.
$ws.PageSetup.Orientation = 1

# Update the view: scroll/selection state.
$ws.Range("C22").Select()
